{"js": "// Replace old equations with new equations, matching the diff exactly.\n// Each old value appears exactly once in the document, so a plain search\n// + replace-insert is unambiguous for every pair.\nconst replacements = [\n  [\"145\u00d76=\", \"560\u00d76=\"],\n  [\"151\u00d79=\", \"597\u00d77=\"],\n  [\"427\u00d74=\", \"361\u00d72=\"],\n  [\"777\u00d72=\", \"160\u00d77=\"],\n  [\"317\u00d77=\", \"250\u00d72=\"],\n  [\"925\u00d74=\", \"107\u00d79=\"],\n  [\"900\u00d77=\", \"234\u00d74=\"],\n  [\"483\u00d73=\", \"405\u00d73=\"],\n  [\"621\u00d74=\", \"699\u00d73=\"],\n  [\"881\u00d77=\", \"515\u00d79=\"],\n  [\"208\u00d75=\", \"401\u00d79=\"],\n  [\"922\u00d74=\", \"728\u00d74=\"],\n  [\"376\u00d78=\", \"323\u00d73=\"],\n  [\"208\u00d73=\", \"835\u00d79=\"],\n  [\"976\u00d76=\", \"303\u00d73=\"],\n  [\"403\u00d78=\", \"801\u00d76=\"],\n  [\"312\u00d78=\", \"640\u00d78=\"],\n  [\"995\u00d73=\", \"311\u00d78=\"],\n  [\"955\u00d76=\", \"978\u00d72=\"],\n  [\"250\u00d76=\", \"404\u00d72=\"],\n  [\"984\u00d74=\", \"873\u00d76=\"],\n  [\"748\u00d74=\", \"620\u00d78=\"],\n  [\"198\u00d77=\", \"162\u00d74=\"],\n  [\"530\u00d76=\", \"147\u00d78=\"],\n  [\"340\u00d73=\", \"465\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "# Replace old equations with new equations, matching the diff exactly.\n# Each old value appears exactly once in the document, so Find/Replace\n# (wdReplaceOne, restricted to each literal match) is unambiguous per pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"145\u00d76=\", \"560\u00d76=\"),\n    @(\"151\u00d79=\", \"597\u00d77=\"),\n    @(\"427\u00d74=\", \"361\u00d72=\"),\n    @(\"777\u00d72=\", \"160\u00d77=\"),\n    @(\"317\u00d77=\", \"250\u00d72=\"),\n    @(\"925\u00d74=\", \"107\u00d79=\"),\n    @(\"900\u00d77=\", \"234\u00d74=\"),\n    @(\"483\u00d73=\", \"405\u00d73=\"),\n    @(\"621\u00d74=\", \"699\u00d73=\"),\n    @(\"881\u00d77=\", \"515\u00d79=\"),\n    @(\"208\u00d75=\", \"401\u00d79=\"),\n    @(\"922\u00d74=\", \"728\u00d74=\"),\n    @(\"376\u00d78=\", \"323\u00d73=\"),\n    @(\"208\u00d73=\", \"835\u00d79=\"),\n    @(\"976\u00d76=\", \"303\u00d73=\"),\n    @(\"403\u00d78=\", \"801\u00d76=\"),\n    @(\"312\u00d78=\", \"640\u00d78=\"),\n    @(\"995\u00d73=\", \"311\u00d78=\"),\n    @(\"955\u00d76=\", \"978\u00d72=\"),\n    @(\"250\u00d76=\", \"404\u00d72=\"),\n    @(\"984\u00d74=\", \"873\u00d76=\"),\n    @(\"748\u00d74=\", \"620\u00d78=\"),\n    @(\"198\u00d77=\", \"162\u00d74=\"),\n    @(\"530\u00d76=\", \"147\u00d78=\"),\n    @(\"340\u00d73=\", \"465\u00d77=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n}\n"}
